$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, since several values (e.g. "1.003", "1.001")
# would otherwise be auto-interpreted by Excel as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.897.17'
$ws.Range("E2").Value = '  +1.43%  '
$ws.Range("D3").Value = '1.664.19'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").Value = '326.90'
$ws.Range("E5").Value = '  +5.98%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("D7").Value = '0.3643'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '47.65'
$ws.Range("E8").Value = '  +1.60%  '
$ws.Range("D9").Value = '0.3271'
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").Value = '1.135'
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("D11").Value = '0.07091'
$ws.Range("E11").Value = '  +1.73%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").Value = '6.073'
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("D14").Value = '19.56'
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").Value = '1.658.84'
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").Value = '6.621'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '0.00001048'
$ws.Range("E17").Value = '  -0.36%  '
$ws.Range("D18").Value = '0.06635'
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.62%  '
$ws.Range("D20").Value = '78.88'
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '5.929'
$ws.Range("E21").Value = '  -0.63%  '
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").Value = '15.82'
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").Value = '12.64'
$ws.Range("E23").Value = '  +4.89%  '
$ws.Range("D24").Value = '24.865.83'
$ws.Range("E24").Value = '  +1.72%  '
$ws.Range("D25").Value = '2.462'
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("D26").Value = '2.427'
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("D27").Value = '149.57'
$ws.Range("E27").Value = '  +1.69%  '
$ws.Range("D28").Value = '18.69'
$ws.Range("E28").Value = '  -1.64%  '
$ws.Range("D29").Value = '1.843.71'
$ws.Range("E29").Value = '  +0.94%  '
$ws.Range("D30").Value = '126.00'
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Value = '1.181'
$ws.Range("E31").Value = '  +9.54%  '
$ws.Range("D32").Value = '4.077'
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("D33").Value = '5.685'
$ws.Range("E33").Value = '  -1.45%  '
$ws.Range("D34").Value = '0.08461'
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").Value = '1.655'
$ws.Range("E35").Value = '  -2.52%  '
$ws.Range("D36").Value = '12.20'
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("D37").Value = '1.291'
$ws.Range("E37").Value = '  +6.56%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.06210'
$ws.Range("E38").Value = '  +1.54%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '5.171'
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.02275'
$ws.Range("E40").Value = '  +1.29%  '
$ws.Range("D41").Value = '0.2080'
$ws.Range("E41").Value = '  +1.02%  '
$ws.Range("D42").Value = '8.277'
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("D44").Value = '0.5945'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").Value = '13.57'
$ws.Range("E45").Value = '  +7.29%  '
$ws.Range("D46").Value = '3.841'
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("D47").Value = '0.5649'
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").Value = '125.54'
$ws.Range("E48").Value = '  +3.18%  '
$ws.Range("D49").Value = '1.954'
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").Value = '0.06986'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").Value = '1.192'
$ws.Range("E51").Value = '  +2.98%  '

# Restore the default (unstyled) cell style so no new formatting is introduced,
# matching the original workbook which had no explicit style on these data cells.
$ws.Range("D2:D51").Style = "Normal"
